$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text: remove trailing "[]" from "PastTimes[]" and "Skills[]"
$ws.Range("I1").Value = "PastTimes"
$ws.Range("J1").Value = "Skills"

# New decimal-string values for column G (Agility), rows 2-14
$agilityValues = @("1.5", "2.5", "12.0", "15.3", "3.2", "6.6", "4.5", "7.7", "4.7", "9.5", "8.6", "4.2", "3.3")

for ($i = 0; $i -lt $agilityValues.Length; $i++) {
    $row = $i + 2
    $ws.Range("G$row").Value = $agilityValues[$i]
}

# Update the active selection to match the new selected cell
$ws.Range("H18").Select()
